# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (cloned from the "2021-Q4" sheet, which has
# the same fund-holding table layout) right after "2021-Q4", and rebuilds the
# "总计" (grand-total) sheet with a new leading row for 2022-Q1 followed by
# the previously-existing rows (shifted down by one).
#
# Sheet reference caveat: worksheet handles captured via Worksheets.Item(...)
# resolve positionally, so any handle obtained *before* a sheet is
# inserted/deleted can silently point at a different sheet afterwards. We
# always (re-)fetch "总计" right before using it, and otherwise keep a single
# stable reference ($src) to "2021-Q4", which never moves.

$wb = $excel.ActiveWorkbook

# --- 1. Delete the old "总计" sheet first -----------------------------------
# It currently holds the highest sheetId (6) in the workbook, so removing it
# now means the next two sheets we create will be assigned sheetId 6 and 7,
# matching "2022-Q1"/"总计" in the target workbook exactly.
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# --- 2. Create "2022-Q1" by cloning "2021-Q4" -------------------------------
# 2021-Q4 already has the right header row (基金代码/基金名称/基金规模/...)
# and cell styles, so copying it preserves sheetPr/pageMargins/styles exactly
# and we only need to overwrite the fund's row-2 values.
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $src)
$newQ1 = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ1.Name = "2022-Q1"

# New fund snapshot for 2022-Q1 (same fund, updated figures). These columns
# are stored as text in the source sheets (not numbers), so assign them with
# a leading apostrophe to force text, then drop the resulting NumberFormat
# style tweak so the cell style stays the plain (unstyled) one, like the rest
# of the data row.
$newQ1.Range("D2").Value = "'54.05"
$newQ1.Range("E2").Value = "'93.05"
$newQ1.Range("F2").Value = "'1.60"
$newQ1.Range("G2").Value = "'0.8648"
$newQ1.Range("D2:G2").ClearFormats()

# --- 3. Create the new "总计" sheet -----------------------------------------
# Clone the freshly-made "2022-Q1" sheet (so sheetPr/pageMargins/etc. match
# the family of sheets) and wipe its cell contents, then rebuild the
# date/count/value table.
$newQ1.Copy($null, $newQ1)
$newTotal = $wb.Worksheets.Item("2022-Q1 (2)")
$newTotal.Name = "总计"
$newTotal.Cells.Clear()

# Header row, re-using the bold/bordered header style (s=2) from the donor
# sheet's header cells.
$src.Range("B1:D1").Copy()
$newTotal.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

# Column A (row index) style, re-using the same style (s=2) applied down
# through row 7.
$src.Range("A2").Copy()
$newTotal.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

# Row 2: new 2022-Q1 entry.
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0.86

# Rows 3-7: previous rows, shifted down by one.
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 1
$newTotal.Range("D3").Value = 0.88

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 1
$newTotal.Range("D4").Value = 0.64

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 1
$newTotal.Range("D5").Value = 0.39

$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 2
$newTotal.Range("D6").Value = 0.24

$newTotal.Range("A7").Value = 5
$newTotal.Range("B7").Value = "2020-Q4"
$newTotal.Range("C7").Value = 1
$newTotal.Range("D7").Value = 0.23
